$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (D) text, new Volume(1h) text (E).
# $null means "leave that column untouched for this row".
$updates = @(
    @{ Row = 2; D = "70.733.97"; E = "  +2.52%  " }
    @{ Row = 3; D = "3.584.10"; E = "  +1.87%  " }
    @{ Row = 4; D = "1.00"; E = $null }
    @{ Row = 5; D = "597.78"; E = "  +1.73%  " }
    @{ Row = 6; D = "173.40"; E = "  +1.65%  " }
    @{ Row = 7; D = "3.577.61"; E = "  +1.86%  " }
    @{ Row = 8; D = "0.617"; E = "  +1.02%  " }
    @{ Row = 9; D = $null; E = "  +0.01%  " }
    @{ Row = 10; D = $null; E = "  +6.06%  " }
    @{ Row = 11; D = "7.36"; E = "  +8.01%  " }
    @{ Row = 12; D = $null; E = "  +1.91%  " }
    @{ Row = 13; D = "46.75"; E = "  -0.42%  " }
    @{ Row = 14; D = $null; E = "  +1.06%  " }
    @{ Row = 15; D = "4.163.28"; E = "  +2.07%  " }
    @{ Row = 16; D = "8.41"; E = "  -0.28%  " }
    @{ Row = 17; D = "614.37"; E = "  -0.23%  " }
    @{ Row = 18; D = "3.578.88"; E = "  +1.59%  " }
    @{ Row = 19; D = "70.798.12"; E = "  +2.47%  " }
    @{ Row = 20; D = $null; E = "  -0.89%  " }
    @{ Row = 21; D = $null; E = "  +0.79%  " }
    @{ Row = 22; D = "0.886"; E = "  +0.46%  " }
    @{ Row = 23; D = "9.27"; E = "  -16.43%  " }
    @{ Row = 24; D = "15.90"; E = "  +0.84%  " }
    @{ Row = 25; D = "97.08"; E = "  +0.48%  " }
    @{ Row = 26; D = "3.77"; E = "  -2.00%  " }
    @{ Row = 27; D = $null; E = "  -0.04%  " }
    @{ Row = 28; D = "2.64"; E = "  +0.98%  " }
    @{ Row = 29; D = "33.97"; E = "  +4.14%  " }
    @{ Row = 30; D = "9.20"; E = "  -0.08%  " }
    @{ Row = 31; D = "8.39"; E = "  -0.86%  " }
    @{ Row = 32; D = "3.07"; E = "  -1.47%  " }
    @{ Row = 33; D = "7.20"; E = "  +4.43%  " }
    @{ Row = 34; D = "649.18"; E = "  +3.47%  " }
    @{ Row = 35; D = $null; E = "  -1.45%  " }
    @{ Row = 36; D = $null; E = "  +6.41%  " }
    @{ Row = 37; D = $null; E = "  -0.80%  " }
    @{ Row = 38; D = "10.84"; E = "  +1.03%  " }
    @{ Row = 39; D = "0.0480"; E = "  +7.29%  " }
    @{ Row = 40; D = $null; E = "  -0.01%  " }
    @{ Row = 41; D = $null; E = "  -0.03%  " }
    @{ Row = 42; D = $null; E = "  +5.59%  " }
    @{ Row = 43; D = "3.394.20"; E = "  +0.79%  " }
    @{ Row = 44; D = $null; E = "  -0.62%  " }
    @{ Row = 45; D = "0.0₃0714"; E = "  +2.99%  " }
    @{ Row = 46; D = "32.96"; E = "  +0.90%  " }
    @{ Row = 47; D = "2.96"; E = "  +6.30%  " }
    @{ Row = 48; D = $null; E = "  +5.07%  " }
    @{ Row = 49; D = $null; E = "  +0.91%  " }
    @{ Row = 50; D = "133.05"; E = "  -0.05%  " }
    @{ Row = 51; D = $null; E = "  -0.13%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Columns D hold price strings that often look like plain numbers
        # (e.g. "1.00", "0.0480"); Excel would silently coerce those to
        # numeric values and drop formatting such as trailing zeros. Force
        # the cell to text first, write the value, then restore the default
        # "Normal" style so no visible formatting change is introduced.
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.NumberFormat = "General"
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
